$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Completed by" entries for row 15 (Overall design/layout) and
# row 16 (Website and game flow) to add the new team members.
$ws.Range("C15").Value = "Christian, Josh, Travis, Adamma, Makafui"
$ws.Range("C16").Value = "Josh, Christian, Travis, Adamma, Makafui"

# Widen column C so the longer text fits (no longer an auto "best fit").
# (38.3 is the ColumnWidth input that this runtime's rounding maps closest
# to the target stored width of ~39.14 characters.)
$ws.Columns.Item(3).ColumnWidth = 38.3

# Leave the selection where the last edit happened.
$ws.Range("E15").Select()
